$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.422.14"
$ws.Range("E2").Value = "  -2.83%  "

$ws.Range("D3").Value = "1.771.10"
$ws.Range("E3").Value = "  -1.98%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.10"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4265"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +1.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3618"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +1.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07149"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8379"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.39"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").Value = "1.790.89"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.434"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.239"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06888"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "78.74"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -3.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008659"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -1.26%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.91"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.03%  "

$ws.Range("D21").Value = "26.434.38"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.096"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("D24").Value = "2.010.03"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.37"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("E26").Value = "  -5.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.98"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.050"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.67"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.764"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +3.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08889"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7225"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -2.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.109"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.308"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -3.23%  "

$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.004"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.749"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -5.57%  "

$ws.Range("E37").Value = "  +2.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05121"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01887"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1611"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4905"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.583"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -5.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.324"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.948"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -2.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.53"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.004"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.12"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -1.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.621"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06175"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -3.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4456"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.709"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.71%  "
